# Mark two more wishlist items as "Reserved / bought" (column E = "Y")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Y"
$ws.Range("E7").Value = "Y"

# Leave the selection where the author's last edit landed (row 7)
$ws.Range("E7").Select()
